$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New shared-string texts (uniqueCount 37 -> 44) ----
$s37 = ' [CS:N]Murkrow[CR]\''s been over there\nfor a while now.'
$s38 = ' It looks like she\''s having a hard\ntime deciding what to order.'
$s39 = 'SCRIPT/P01P04A/um1102.ssb'
$s40 = ' [CS:N]Маркроу[CR] здесь уже долго стоит.'
$s41 = ' Похоже, она не может решить,\nчто ей заказать.'
$s42 = ' [CS:N]Íàñëñïô[CR] èäåòû ôçå äïìãï òóïéó.'
$s43 = ' Ðïöïçå, ïîà îå íïçåó ñåšéóû,\nœóï åê èàëàèàóû.'

# ---- Row 11: mirror the formatting of row 2 (A/B = style "4", C/D/E = style "5") ----
$ws.Range("A2:E2").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)

# ---- Row 12: mirror the formatting of row 5 (no A cell; B = style "4", C/D/E = style "5") ----
$ws.Range("B5:E5").Copy()
$ws.Range("B12:E12").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Values: written in the same order the shared-string table expects
#      them (English column first, then the filename, then Russian, then
#      the "converted"/ciphered column) so the appended <si> entries line
#      up with indices 37-43 exactly as in the target workbook. ----
$ws.Range("C11").Value = $s37
$ws.Range("C12").Value = $s38
$ws.Range("A11").Value = $s39
$ws.Range("D11").Value = $s40
$ws.Range("D12").Value = $s41
$ws.Range("E11").Value = $s42
$ws.Range("E12").Value = $s43

$ws.Range("B11").Value = 307
$ws.Range("B12").Value = 316

# ---- Row heights (wrap-text auto-fit sizes for the new content) ----
$ws.Rows.Item(11).RowHeight = 43.2
$ws.Rows.Item(12).RowHeight = 21.6

# ---- View state: move selection to the newly added cell ----
$ws.Range("C12").Select()
